$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(87, 3).Value2 = 5.67151440450491
$ws.Cells.Item(87, 5).Value2 = 8.51051752411876
$ws.Cells.Item(87, 8).Value2 = 7.81619007238754
$ws.Cells.Item(88, 3).Value2 = 2.08818981574577
$ws.Cells.Item(88, 5).Value2 = -0.0424822104396039
$ws.Cells.Item(88, 8).Value2 = -1.18899392005855
$ws.Cells.Item(89, 3).Value2 = 27.6327956778756
$ws.Cells.Item(89, 5).Value2 = 105.514216332317
$ws.Cells.Item(89, 8).Value2 = 104.63304871606
$ws.Cells.Item(90, 3).Value2 = 52.5194921267955
$ws.Cells.Item(90, 5).Value2 = 96.0957168611856
$ws.Cells.Item(90, 8).Value2 = 95.9179395174032
$ws.Cells.Item(91, 3).Value2 = 31.0421778300851
$ws.Cells.Item(91, 5).Value2 = -77.398739662723
$ws.Cells.Item(91, 8).Value2 = -77.5372853184251
$ws.Cells.Item(92, 3).Value2 = 12.7499706683285
$ws.Cells.Item(92, 5).Value2 = -73.2113108574661
$ws.Cells.Item(92, 8).Value2 = -73.4991182613114
$ws.Cells.Item(93, 3).Value2 = -17.9302595466106
$ws.Cells.Item(93, 5).Value2 = -17.206704527439
$ws.Cells.Item(93, 8).Value2 = -17.4556415051167
$ws.Cells.Item(94, 3).Value2 = -45.7049371337172
$ws.Cells.Item(94, 5).Value2 = -15.0029934872405
$ws.Cells.Item(94, 8).Value2 = -15.0208298620717
$ws.Cells.Item(95, 3).Value2 = -30.3476123982645
$ws.Cells.Item(95, 5).Value2 = -15.9694407209125
$ws.Cells.Item(95, 8).Value2 = -15.9365776351469
$ws.Cells.Item(96, 3).Value2 = -15.8999659597439
$ws.Cells.Item(96, 5).Value2 = -15.4207251033835
$ws.Cells.Item(96, 8).Value2 = -15.4053228961234
$ws.Cells.Item(97, 3).Value2 = -12.0302035535194
$ws.Cells.Item(97, 5).Value2 = -1.72765490254099
$ws.Cells.Item(97, 8).Value2 = -1.73189732289647
$ws.Cells.Item(98, 3).Value2 = -8.77240957312753
$ws.Cells.Item(98, 5).Value2 = -1.97181756567307
$ws.Cells.Item(98, 8).Value2 = -1.94439192954412
$ws.Cells.Item(99, 3).Value2 = -4.99458254316726
$ws.Cells.Item(99, 5).Value2 = -0.858132601071416
$ws.Cells.Item(99, 8).Value2 = -0.845440452562746
$ws.Cells.Item(100, 3).Value2 = -2.9898961633646
$ws.Cells.Item(100, 5).Value2 = -7.40197958417291
$ws.Cells.Item(100, 8).Value2 = -7.39457014889271
$ws.Cells.Item(101, 3).Value2 = -4.37122482421409
$ws.Cells.Item(101, 5).Value2 = -7.25296954593898
$ws.Cells.Item(101, 8).Value2 = -7.2153281723287
$ws.Cells.Item(102, 3).Value2 = -5.54307374294341
$ws.Cells.Item(102, 5).Value2 = -6.65921324059034
$ws.Cells.Item(102, 8).Value2 = -6.57623873829117
$ws.Cells.Item(103, 3).Value2 = -6.94902555639709
$ws.Cells.Item(103, 5).Value2 = -6.48193985488613
$ws.Cells.Item(103, 8).Value2 = -6.39688849011566
$ws.Cells.Item(104, 3).Value2 = -5.12107276508766
$ws.Cells.Item(104, 5).Value2 = -0.0901684189351754
$ws.Cells.Item(104, 8).Value2 = -0.0111309079239056
$ws.Cells.Item(105, 3).Value2 = -3.30607380498818
$ws.Cells.Item(105, 5).Value2 = 0.00702629445892557
$ws.Cells.Item(105, 8).Value2 = 0.0955870571189224
$ws.Cells.Item(106, 3).Value2 = -1.67610129381237
$ws.Cells.Item(106, 5).Value2 = -0.139323195887082
$ws.Cells.Item(106, 8).Value2 = -0.037591231689661
$ws.Cells.Item(107, 3).Value2 = -0.0904471290626034
$ws.Cells.Item(107, 8).Value2 = -0.0974341700321491
$ws.Cells.Item(108, 3).Value2 = -0.10273582330058
$ws.Cells.Item(108, 8).Value2 = -0.191292248379047
$ws.Cells.Item(109, 3).Value2 = -0.139323195887082
$ws.Cells.Item(109, 8).Value2 = -0.222056238134473
$ws.Cells.Item(110, 3).Value2 = -0.139323195887082
$ws.Cells.Item(111, 3).Value2 = -0.139323195887082
$ws.Cells.Item(112, 3).Value2 = -0.139323195887082
$ws.Cells.Item(113, 3).Value2 = -0.139323195887082
$ws.Cells.Item(114, 3).Value2 = -0.139323195887082
$ws.Cells.Item(115, 3).Value2 = -0.139323195887082
$ws.Cells.Item(116, 3).Value2 = -0.139323195887082
$ws.Cells.Item(117, 3).Value2 = -0.139323195887082
$ws.Cells.Item(118, 3).Value2 = -0.139323195887082
$ws.Cells.Item(119, 3).Value2 = -0.139323195887082
$ws.Cells.Item(120, 3).Value2 = -0.139323195887082
$ws.Cells.Item(121, 3).Value2 = -0.139323195887082
$ws.Cells.Item(122, 3).Value2 = -0.139323195887082
$ws.Cells.Item(123, 3).Value2 = -0.139323195887082
$ws.Cells.Item(124, 3).Value2 = -0.139323195887082
$ws.Cells.Item(125, 3).Value2 = -0.139323195887082
$ws.Cells.Item(126, 3).Value2 = -0.139323195887082
$ws.Cells.Item(127, 3).Value2 = -0.139323195887082
$ws.Cells.Item(128, 3).Value2 = -0.139323195887082
$ws.Cells.Item(129, 3).Value2 = -0.139323195887082
$ws.Cells.Item(130, 3).Value2 = -0.139323195887082
